# Trade #92 (EMAArbitrage) closes early at 2026-02-18 00:23:47, and a new
# Trade #121 (MarketMaking) opens at 2026-02-18 00:23:41 - update the
# Summary, Strategy Status, All Trades, MarketMaking and EMAArbitrage
# sheets to reflect both events.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: bump total trade count, recompute win rate %
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 92
$summary.Range("B9").Value = 48.91

# ---------------------------------------------------------------------
# Strategy Status sheet: EMAArbitrage is row 2
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D2").Value = 6
$status.Range("G2").Value = 50

# ---------------------------------------------------------------------
# All Trades sheet: close out trade #92 (row 93) and append the new
# trade #121 (row 122)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(93, 7).Value = 0.98
$allTrades.Cells.Item(93, 8).Value = "CLOSED"
$allTrades.Cells.Item(93, 11).Value = 100.27
$allTrades.Cells.Item(93, 12).Value = "early_exit"
$allTrades.Cells.Item(93, 13).Value = 0.13

$allTrades.Cells.Item(122, 1).Value = 121
$c = $allTrades.Cells.Item(122, 2)
$c.NumberFormat = "@"
$c.Value = "2026-02-18"
$allTrades.Cells.Item(122, 3).Value = "00:23:41"
$allTrades.Cells.Item(122, 4).Value = "MarketMaking"
$allTrades.Cells.Item(122, 5).Value = "DOWN"
$allTrades.Cells.Item(122, 6).Value = 0.98
$allTrades.Cells.Item(122, 8).Value = "OPEN"
$allTrades.Cells.Item(122, 9).Value = 0
$allTrades.Cells.Item(122, 10).Value = 0
$allTrades.Cells.Item(122, 11).Value = 99.410254715139
$allTrades.Cells.Item(122, 13).Value = 0
$allTrades.Cells.Item(122, 14).Value = 0
$allTrades.Cells.Item(122, 15).Value = 0
$allTrades.Cells.Item(122, 16).Value = 0.6
$allTrades.Cells.Item(122, 17).Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet: append the new trade #121 as row 42
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Cells.Item(42, 1).Value = 121
$c2 = $mm.Cells.Item(42, 2)
$c2.NumberFormat = "@"
$c2.Value = "2026-02-18"
$mm.Cells.Item(42, 3).Value = "00:23:41"
$mm.Cells.Item(42, 4).Value = "MarketMaking"
$mm.Cells.Item(42, 5).Value = "DOWN"
$mm.Cells.Item(42, 6).Value = 0.98
$mm.Cells.Item(42, 8).Value = "OPEN"
$mm.Cells.Item(42, 9).Value = 0
$mm.Cells.Item(42, 10).Value = 0
$mm.Cells.Item(42, 11).Value = 99.410254715139
$mm.Cells.Item(42, 12).Value = 0
$mm.Cells.Item(42, 13).Value = 0
$mm.Cells.Item(42, 14).Value = 0.6
$mm.Cells.Item(42, 15).Value = "Normal spread capture: 198 bps"
$mm.Cells.Item(42, 17).Value = 0

# ---------------------------------------------------------------------
# EMAArbitrage sheet: close out trade #92, which is row 7 here
# ---------------------------------------------------------------------
$ema = $wb.Worksheets.Item("EMAArbitrage")
$ema.Cells.Item(7, 7).Value = 0.98
$ema.Cells.Item(7, 8).Value = "CLOSED"
$ema.Cells.Item(7, 11).Value = 100.27
$ema.Cells.Item(7, 16).Value = "early_exit"
$ema.Cells.Item(7, 17).Value = 0.13
